$wb = $excel.ActiveWorkbook

# --- Sheet1: "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Update Version value
$meta.Range("B3").Value = "6.0.0"

# Update Date value
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Fill in Publisher value
$meta.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact" row (row 11); everything below shifts up.
$meta.Rows.Item(11).Delete()

# Row 10 (previously the first "Contact" row) becomes "Jurisdiction"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet2: "Elements" ---
$elem = $wb.Worksheets.Item("Elements")

# Row 2 (Path = Extension) gets a real Short/Definition update
$elem.Range("K2").Value = "Claim Snapshot Provider Zip Code"
$elem.Range("L2").Value = "Original provider postal code, as reported on the claim"
